# IPASS Logboek update: add "Woensdag" week-3 log entries (rows 17-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new log entries. Cells are written in the same order the
# corresponding strings first appear in the target shared-string table so
# that the generated xl/sharedStrings.xml ordering matches exactly.
$ws.Range("A18").Value = "acknowedge bytes"
$ws.Range("B17").Value = "major restructuring library"
$ws.Range("B18").Value = "dummy pins"
$ws.Range("B20").Value = "all_from_pin_out_t"
$ws.Range("B21").Value = "acknowledge bytes weggehaald"
$ws.Range("B22").Value = "due_remote_primitives"
$ws.Range("B23").Value = "formatting opmooing"
$ws.Range("B19").Value = "hc595 bitbanged spi protocol"
$ws.Range("C17").Value = "removed hc595"
$ws.Range("C20").Value = "port_out_from_pins_t"
$ws.Range("C19").Value = "port_in_out_from_pins_t"
$ws.Range("C21").Value = "port_in_from_pins_t"
$ws.Range("C22").Value = "all_from_port_out_t"
$ws.Range("C18").Value = "port_out primitives"
$ws.Range("C23").Value = "all_from_pin_out_t tested"
$ws.Range("C24").Value = "port_out_from_pins_t tested"
$ws.Range("C25").Value = "port_in_out_from_pins_t tested"

# Widen column C slightly to fit the new, longer entries.
$ws.Columns.Item(3).ColumnWidth = 30.5

# Leave the selection where the author left it after typing the last entry.
$ws.Range("B28").Select()
